$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 33695.279999999999
    3  = 23068.243999999999
    4  = 29389.78
    5  = 33172.695
    6  = 30956.425999999999
    7  = 30220.115000000002
    8  = 33588.785000000003
    9  = 43545.065999999999
    10 = 40279.163999999997
    11 = 46516.71
    12 = 55686.1
    13 = 68610.05
    14 = 60287.964999999997
    15 = 63304.714999999997
    16 = 68542.880000000005
    17 = 76550.31
    18 = 82388.149999999994
    19 = 99322.559999999998
    20 = 87800.39
    21 = 94237.9
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Update the worksheet view: scroll and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("F35").Select()
